# Auto-generated script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.834.09"
$ws.Range("E2").Value = "  -4.06%  "
$ws.Range("D3").Value = "3.587.80"
$ws.Range("E3").Value = "  -4.35%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.73"
$ws.Range("E5").Value = "  -4.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.62"
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.602"
$ws.Range("E7").Value = "  -5.68%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.666"
$ws.Range("E9").Value = "  -8.56%  "
$ws.Range("E10").Value = "  -14.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.89"
$ws.Range("E11").Value = "  -7.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000243"
$ws.Range("E12").Value = "  -18.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.82"
$ws.Range("E13").Value = "  -7.69%  "
$ws.Range("D14").Value = "4.169.87"
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("D15").Value = "3.589.99"
$ws.Range("E15").Value = "  -3.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.126"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "66.629.34"
$ws.Range("E17").Value = "  -3.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.19"
$ws.Range("E18").Value = "  -7.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.07"
$ws.Range("E19").Value = "  -7.55%  "
$ws.Range("E20").Value = "  -7.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.14"
$ws.Range("E21").Value = "  -6.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("E22").Value = "  -9.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.28"
$ws.Range("E23").Value = "  -6.56%  "
$ws.Range("E24").Value = "  -8.71%  "
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("E26").Value = "  -7.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.23"
$ws.Range("E27").Value = "  -8.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.55"
$ws.Range("E28").Value = "  -10.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.81"
$ws.Range("E29").Value = "  -8.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.77"
$ws.Range("E30").Value = "  -7.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.68"
$ws.Range("E31").Value = "  -10.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "65.64"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.73"
$ws.Range("E33").Value = "  -7.34%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "579.26"
$ws.Range("E34").Value = "  -5.49%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.110"
$ws.Range("E35").Value = "  -7.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.96"
$ws.Range("E36").Value = "  -8.02%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.368"
$ws.Range("E39").Value = "  -9.71%  "
$ws.Range("E40").Value = "  -5.31%  "
$ws.Range("D41").Value = "0.0₃0718"
$ws.Range("E41").Value = "  -21.26%  "
$ws.Range("E42").Value = "  -11.13%  "
$ws.Range("E43").Value = "  -9.10%  "
$ws.Range("E44").Value = "  -4.55%  "
$ws.Range("D45").Value = "2.660.26"
$ws.Range("E45").Value = "  -5.34%  "
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "139.75"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  -17.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.34"
$ws.Range("E50").Value = "  -11.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.54"
$ws.Range("E51").Value = "  -9.90%  "
